$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Telavi")

# Update the row 4 values (E4:K4) with the new figures
$ws.Range("E4").Value = 2299
$ws.Range("F4").Value = 2304
$ws.Range("G4").Value = 2295
$ws.Range("H4").Value = 2270
$ws.Range("I4").Value = 2336
$ws.Range("J4").Value = 2387
$ws.Range("K4").Value = 2447

# Update the active selection from A1:K1 to A3
$ws.Range("A3").Select()
